# Revert the cached "Update automatically" footer date field text from
# 2/14/2025 back to 2/3/2025 across the slide master, every slide layout,
# and the notes master (the only user-visible content change in the diff).

$p = $ppt.ActivePresentation

$oldText = "2/14/2025"
$newText = "2/3/2025"

function Update-DateShapes($shapeRange) {
    for ($i = 1; $i -le $shapeRange.Count; $i++) {
        $sh = $shapeRange.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $tf = $sh.TextFrame
            if ($tf.HasText -eq -1) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldText) {
                    $tr.Text = $newText
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
